$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ I = 0.7832871899279559; N = 1.177962427763276 }
    3  = @{ I = 1.539959621819876;  N = 1.283842295398743 }
    4  = @{ I = 0.9944570360131537; N = 1.233158504608332 }
    5  = @{ I = 0.3243737962186763; N = 1.050533268808214 }
    6  = @{ I = 0.5109532906258769; N = 1.038288204613594 }
    7  = @{ I = 0.2998897651543369; N = 1.094407364347728 }
    8  = @{ I = 0.8623813251500541; N = 1.03543469952145 }
    9  = @{ I = 1.506355255511725;  N = 0.9958432720815682 }
    10 = @{ I = 0.2796576970884782; N = 1.154371560340464 }
    11 = @{ I = 2.308041632221933;  N = 0.9663359963356604 }
    12 = @{ I = 3.808791409918911;  N = 1.414455928493327 }
    13 = @{ I = 4.563680258747716;  N = 1.782199216675866 }
    14 = @{ I = 1.181503973411165;  N = 0.7099136471736143 }
    15 = @{ I = 0.9183962002005244; N = 0.7362457997114868 }
    16 = @{ I = 0.274194772300126;  N = 0.9130761450487295 }
    17 = @{ I = 0.2706642960501706; N = 0.9086830241910144 }
    18 = @{ I = 0.4927394096119291; N = 0.7965212918588809 }
    19 = @{ I = 1.087654735259869;  N = 0.721813620111174 }
    20 = @{ I = 0.5599703205083473; N = 0.78492167307603 }
    21 = @{ I = 0.3592403072165918; N = 0.8303565284176275 }
    22 = @{ I = 0.2842370595323462; N = 0.8605421488380638 }
    23 = @{ I = 0.2733077081458004; N = 0.8662700479028254 }
    24 = @{ I = 1.448329008321893;  N = 0.7044970904417205 }
    25 = @{ I = 0.5055471945442691; N = 0.7953108772398096 }
}

foreach ($row in $values.Keys) {
    $ws.Range("I$row").Value = $values[$row].I
    $ws.Range("N$row").Value = $values[$row].N
}
